# "derniere modif avant la fin"
#
# 1) Insert a new "Conclusion" slide right before the final slide
#    (the closing "C'est deja finis !" slide gets pushed one position later).
# 2) Fill in the Conclusion slide's title + content placeholder text.
# 3) Fix the "Booky" -> "Booki" typo on slide 3.
# 4) Refresh the cached date placeholder text (28/10/2021 -> 01/11/2021)
#    on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- 1) Insert the new slide before the last one, using the same
#        "Title and Content" layout as the rest of the deck. ---
$lastIndex = $p.Slides.Count
$newSlide = $p.Slides.Add($lastIndex, 2)

# --- 2) Title placeholder ---
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Conclusion"

# --- Content placeholder ---
$body = $newSlide.Shapes.Item(2)
$paragraphs = @(
    "Aime et pas aimer ",
    "Expliquer un peu plus des elements de details",
    "Retour d" + [char]0x2019 + "eperience",
    "Difficulte pour definir la precision a cause de maquette pdf et psd",
    "Diff pour integration mobile sur maquette iphone 8 alors que sur brief 768px",
    "Parler du rendu populaire a cause de la width de la descritpion"
)
$body.TextFrame.TextRange.Text = [string]::Join([char]13, $paragraphs)

# --- 3) Typo fix: "Booky" -> "Booki" on slide 3 ---
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shp = $s3.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -match "Booky") {
                [void]$tr.Replace("Booky", "Booki")
            }
        }
    }
}

# --- 4) Refresh the cached "today" date text across master + layouts ---
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "28/10/2021") {
                    $shp.TextFrame.TextRange.Text = "01/11/2021"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}
